$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.23'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.364'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05868'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.387'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.375'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8134'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9904'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01121'
$ws.Range("E10").Value = '9OneONEBestin24h'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1421'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03537'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07353'
$ws.Range("E13").Value = '12MandalaExchangeTokenMDX'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03024'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09398'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001587'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04820'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006298'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004084'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009859'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.00009704'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.687'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.199'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("B24").Value = 'BitpandaEcosystemToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.3253'
$ws.Range("E24").Value = '23BitpandaEcosystemTokenBEST'
$ws.Range("B25").Value = 'ProBitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1296'
$ws.Range("E25").Value = '24ProBitTokenPROB'
$ws.Range("B26").Value = 'MCDex'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.164'
$ws.Range("E26").Value = '25MCDexMCB'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03845'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006483'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1072'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.005760'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005666'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6514'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07890'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'

Write-Host "Applied all cell updates"
